$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2 (income / NULL MODEL)
$ws1.Range("C2").Value = 0.97278525902507
$ws1.Range("D2").Value = 0.0272435083185367
$ws1.Range("E2").Value = 1.00002876734361
$ws1.Range("J2").Value = 0.0272427246177168
$ws1.Range("K2").Value = 0.0140584110959639
$ws1.Range("L2").Value = 0.169719947404431
$ws1.Range("M2").Value = 0.230056270645448
$ws1.Range("N2").Value = 0.183778358500395

# Row 3 (income / CONDITIONAL MODEL)
$ws1.Range("F3").Value = 0.912447200068311
$ws1.Range("G3").Value = 0.0140588155191065

# Row 4 (income / COMPLETE MODEL)
$ws1.Range("H4").Value = 0.742722370271836
$ws1.Range("I4").Value = 0
$ws1.Range("O4").Value = 0.257298995263165

# --- Sheet 2: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2 (Sibcorr / income)
$ws2.Range("C2").Value = 0.0272427246177168
$ws2.Range("D2").Value = -0.031019051702992
$ws2.Range("E2").Value = 0.0855045009384256

# Row 3 (IOLIB / income)
$ws2.Range("C3").Value = 0.183778358500395
$ws2.Range("D3").Value = 0.148236280230852
$ws2.Range("E3").Value = 0.219320436769938

# Row 4 (IORAD / income)
$ws2.Range("C4").Value = 0.257298995263165
$ws2.Range("D4").Value = 0.216543929670222
$ws2.Range("E4").Value = 0.298054060856107
